$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.074.69"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").Value = "1.552.04"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'287.08"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").Value = "'0.3841"
$ws.Range("E7").Value = "  +3.25%  "
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("D9").Value = "'43.61"
$ws.Range("E9").Value = "  -9.54%  "
$ws.Range("D10").Value = "'1.125"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").Value = "'0.07361"
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'20.04"
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").Value = "'5.792"
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("D15").Value = "1.582.26"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "'6.752"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "'0.00001086"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").Value = "'0.06625"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("D19").Value = "'85.71"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "'6.376"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").Value = "'16.02"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").Value = "'11.67"
$ws.Range("E23").Value = "  -3.23%  "
$ws.Range("D24").Value = "22.080.09"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("E25").Value = "  -3.73%  "
$ws.Range("D26").Value = "'2.492"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").Value = "'150.38"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").Value = "'19.10"
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("D29").Value = "'4.925"
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("D30").Value = "1.758.24"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").Value = "'121.42"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").Value = "'1.078"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").Value = "'5.882"
$ws.Range("E33").Value = "  -3.82%  "
$ws.Range("D34").Value = "'1.900"
$ws.Range("E34").Value = "  -5.55%  "
$ws.Range("D35").Value = "'0.08226"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "'9.236"
$ws.Range("E36").Value = "  -5.08%  "
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("D38").Value = "'0.02315"
$ws.Range("E38").Value = "  -5.64%  "
$ws.Range("D39").Value = "'5.263"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").Value = "'0.2149"
$ws.Range("E40").Value = "  -5.30%  "
$ws.Range("D41").Value = "'1.232"
$ws.Range("E41").Value = "  -4.16%  "
$ws.Range("D42").Value = "'11.00"
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "'0.6000"
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").Value = "'3.727"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").Value = "'0.5810"
$ws.Range("E47").Value = "  -5.32%  "
$ws.Range("D48").Value = "'1.968"
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("D49").Value = "'121.92"
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("E50").Value = "  -2.96%  "
$ws.Range("D51").Value = "'0.07020"
